$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 14001.75
$ws.Range("J28").Value = 25003.5
$ws.Range("L28").Value = 25003.5
$ws.Range("N28").Value = -25973.5
$ws.Range("H40").Value = 38506.168
$ws.Range("I40").Value = 27759.25
$ws.Range("K40").Value = 27759.25
$ws.Range("M40").Value = -27584.25
$ws.Range("H48").Value = 8999
$ws.Range("H56").Value = 8999
$ws.Range("H62").Value = 4118.8335
$ws.Range("I62").Value = 3603
$ws.Range("K62").Value = 3603
$ws.Range("M62").Value = -2979
$ws.Range("H65").Value = 4118.8335
$ws.Range("I65").Value = 3603
$ws.Range("K65").Value = 18015
$ws.Range("M65").Value = -14895
$ws.Range("H69").Value = 9005.5
$ws.Range("I69").Value = 8013
$ws.Range("K69").Value = 24039
$ws.Range("M69").Value = -23165
$ws.Range("H72").Value = 9005.5
$ws.Range("I72").Value = 8013
$ws.Range("K72").Value = 72117
$ws.Range("M72").Value = -67749
$ws.Range("H113").Value = 35741850
$ws.Range("J113").Value = 58867344
$ws.Range("L113").Value = 58867344
$ws.Range("N113").Value = -58873852
$ws.Range("H127").Value = 2651.5676
$ws.Range("I127").Value = 1141.2
$ws.Range("K127").Value = 3423.6
$ws.Range("M127").Value = 1536.4
$ws.Range("H132").Value = 126872.98
$ws.Range("I132").Value = 412946.9
$ws.Range("K132").Value = 1238840.7
$ws.Range("M132").Value = -1236310.7
$ws.Range("H135").Value = 6296.722
$ws.Range("I135").Value = 495
$ws.Range("K135").Value = 4455
$ws.Range("M135").Value = -1920
$ws.Range("H137").Value = 55561576
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("H138").Value = 5674.3784
$ws.Range("J138").Value = 6752.768
$ws.Range("L138").Value = 20258.304
$ws.Range("N138").Value = -30538.304
$ws.Range("H141").Value = 4815.6
$ws.Range("I141").Value = 4402.615
$ws.Range("K141").Value = 13207.845
$ws.Range("M141").Value = -8027.844999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3382.2
$ws.Range("I45").Value = 2470.6667
$ws.Range("J45").Value = 4749.5
$ws.Range("K45").Value = 2470.6667
$ws.Range("L45").Value = 4749.5
$ws.Range("M45").Value = -2093.6667
$ws.Range("N45").Value = -5503.5
$ws.Range("H61").Value = 9442014
$ws.Range("I61").Value = 12201302
$ws.Range("K61").Value = 12201302
$ws.Range("M61").Value = -12201090
$ws.Range("H74").Value = 26789860
$ws.Range("I74").Value = 62501052
$ws.Range("J74").Value = 6468.6875
$ws.Range("K74").Value = 62501052
$ws.Range("L74").Value = 6468.6875
$ws.Range("M74").Value = -62500178
$ws.Range("N74").Value = -8216.6875
$ws.Range("H77").Value = 26789860
$ws.Range("I77").Value = 62501052
$ws.Range("J77").Value = 6468.6875
$ws.Range("K77").Value = 312505260
$ws.Range("L77").Value = 32343.4375
$ws.Range("M77").Value = -312500892
$ws.Range("N77").Value = -41079.4375
$ws.Range("H102").Value = 966.625
$ws.Range("I102").Value = 997.46155
$ws.Range("J102").Value = 833
$ws.Range("K102").Value = 997.46155
$ws.Range("L102").Value = 833
$ws.Range("M102").Value = 624.53845
$ws.Range("N102").Value = -4077
$ws.Range("H110").Value = 4465.6816
$ws.Range("I110").Value = 4030.889
$ws.Range("J110").Value = 6422.25
$ws.Range("K110").Value = 4030.889
$ws.Range("L110").Value = 6422.25
$ws.Range("M110").Value = -1985.889
$ws.Range("N110").Value = -10512.25
$ws.Range("H122").Value = 27783440
$ws.Range("I122").Value = 41669660
$ws.Range("J122").Value = 11000.5
$ws.Range("K122").Value = 125008980
$ws.Range("L122").Value = 33001.5
$ws.Range("M122").Value = -125006530
$ws.Range("N122").Value = -37901.5
$ws.Range("H136").Value = 9442014
$ws.Range("I136").Value = 12201302
$ws.Range("K136").Value = 36603906
$ws.Range("M136").Value = -36601356

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 47621044
$ws.Range("I86").Value = 1673.5
$ws.Range("K86").Value = 1673.5
$ws.Range("M86").Value = -550.5
$ws.Range("H89").Value = 47621044
$ws.Range("I89").Value = 1673.5
$ws.Range("K89").Value = 8367.5
$ws.Range("M89").Value = -2751.5
$ws.Range("H141").Value = 79805.2
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4213.551
$ws.Range("I31").Value = 2890.375
$ws.Range("K31").Value = 2890.375
$ws.Range("M31").Value = -2595.375
$ws.Range("H34").Value = 4213.551
$ws.Range("I34").Value = 2890.375
$ws.Range("K34").Value = 2890.375
$ws.Range("M34").Value = -2688.375
$ws.Range("H86").Value = 5547.8335
$ws.Range("I86").Value = 5457.4
$ws.Range("K86").Value = 5457.4
$ws.Range("M86").Value = -4334.4
$ws.Range("H89").Value = 5547.8335
$ws.Range("I89").Value = 5457.4
$ws.Range("K89").Value = 27287
$ws.Range("M89").Value = -21671
$ws.Range("H99").Value = 7414.5
$ws.Range("I99").Value = 7455
$ws.Range("J99").Value = 7365
$ws.Range("K99").Value = 7455
$ws.Range("L99").Value = 7365
$ws.Range("M99").Value = -5957
$ws.Range("N99").Value = -10361
$ws.Range("H126").Value = 7414.5
$ws.Range("I126").Value = 7455
$ws.Range("J126").Value = 7365
$ws.Range("K126").Value = 22365
$ws.Range("L126").Value = 22095
$ws.Range("M126").Value = -19895
$ws.Range("N126").Value = -27035
$ws.Range("H132").Value = 2208.3
$ws.Range("I132").Value = 1014.1667
$ws.Range("K132").Value = 3042.5001
$ws.Range("M132").Value = -512.5001000000002
$ws.Range("H134").Value = 4121.2
$ws.Range("I134").Value = 3745
$ws.Range("K134").Value = 11235
$ws.Range("M134").Value = -8700

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 49693.41
$ws.Range("J80").Value = 6090.636
$ws.Range("L80").Value = 6090.636
$ws.Range("N80").Value = -8086.636
$ws.Range("H83").Value = 49693.41
$ws.Range("J83").Value = 6090.636
$ws.Range("L83").Value = 30453.18
$ws.Range("N83").Value = -40437.18
$ws.Range("H102").Value = 12411.5
$ws.Range("I102").Value = 14901.777
$ws.Range("J102").Value = 9921.223
$ws.Range("K102").Value = 14901.777
$ws.Range("L102").Value = 9921.223
$ws.Range("M102").Value = -13279.777
$ws.Range("N102").Value = -13165.223
$ws.Range("H113").Value = 528478.3
$ws.Range("I113").Value = 1184063.8
$ws.Range("J113").Value = 4009.95
$ws.Range("K113").Value = 1184063.8
$ws.Range("L113").Value = 4009.95
$ws.Range("M113").Value = -1181893.8
$ws.Range("N113").Value = -8349.95
$ws.Range("H126").Value = 4458.3335
$ws.Range("I126").Value = 2277.7778
$ws.Range("J126").Value = 11000
$ws.Range("K126").Value = 6833.3334
$ws.Range("L126").Value = 33000
$ws.Range("M126").Value = -4363.3334
$ws.Range("N126").Value = -37940
$ws.Range("H132").Value = 7041.8096
$ws.Range("I132").Value = 6380.2144
$ws.Range("J132").Value = 8365
$ws.Range("K132").Value = 19140.6432
$ws.Range("L132").Value = 25095
$ws.Range("M132").Value = -16610.6432
$ws.Range("N132").Value = -30155
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 23813584
$ws.Range("I40").Value = 3130.625
$ws.Range("J40").Value = 55560856
$ws.Range("K40").Value = 3130.625
$ws.Range("L40").Value = 55560856
$ws.Range("M40").Value = -2994.625
$ws.Range("N40").Value = -55561128
$ws.Range("H46").Value = 5174.643
$ws.Range("J46").Value = 5883.636
$ws.Range("L46").Value = 5883.636
$ws.Range("N46").Value = -6259.636
$ws.Range("H93").Value = 890.0455
$ws.Range("I93").Value = 912.7857
$ws.Range("J93").Value = 850.25
$ws.Range("K93").Value = 912.7857
$ws.Range("L93").Value = 850.25
$ws.Range("M93").Value = 335.2143
$ws.Range("N93").Value = -3346.25
$ws.Range("H136").Value = 5382.2446
$ws.Range("I136").Value = 4980.054
$ws.Range("K136").Value = 14940.162
$ws.Range("M136").Value = -12390.162

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2231000.2
$ws.Range("I5").Value = 10000001
$ws.Range("K5").Value = 10000001
$ws.Range("M5").Value = -9999889
$ws.Range("H74").Value = 17087.25
$ws.Range("I74").Value = 15875
$ws.Range("K74").Value = 15875
$ws.Range("M74").Value = -14939
$ws.Range("H77").Value = 17087.25
$ws.Range("I77").Value = 15875
$ws.Range("K77").Value = 47625
$ws.Range("M77").Value = -42945
$ws.Range("H132").Value = 34745056
$ws.Range("I132").Value = 3993990
$ws.Range("K132").Value = 11981970
$ws.Range("M132").Value = -11979440
$ws.Range("H136").Value = 9101.642
$ws.Range("I136").Value = 3839.9092
$ws.Range("J136").Value = 9928.485
$ws.Range("K136").Value = 11519.7276
$ws.Range("L136").Value = 29785.455
$ws.Range("M136").Value = -8969.7276
$ws.Range("N136").Value = -34885.455
